$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source publishes a new weekly record for
# "Feria Lagunitas de Puerto Montt - Zanahoria". It lands as the newest row,
# right after the header block of existing rows, which pushes every
# following data row down by one (506 data rows -> 507 data rows).
$ws.Rows(406).Insert()

$newRow = 406
$ws.Cells.Item($newRow, 1).Value2  = 4
$ws.Cells.Item($newRow, 2).Value   = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($newRow, 3).Value   = "Los Lagos"
$ws.Cells.Item($newRow, 4).Value2  = 44964
$ws.Cells.Item($newRow, 5).Value2  = 10
$ws.Cells.Item($newRow, 6).Value2  = 100114013
$ws.Cells.Item($newRow, 7).Value   = "Zanahoria"
$ws.Cells.Item($newRow, 8).Value   = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value   = "Primera"
$ws.Cells.Item($newRow, 10).Value2 = 700
$ws.Cells.Item($newRow, 11).Value2 = 12000
$ws.Cells.Item($newRow, 12).Value2 = 13000
$ws.Cells.Item($newRow, 13).Value2 = 12500
$ws.Cells.Item($newRow, 14).Value  = "$/saco 20 kilos"
$ws.Cells.Item($newRow, 15).Value  = "Región Metropolitana"
$ws.Cells.Item($newRow, 16).Value2 = 625
$ws.Cells.Item($newRow, 17).Value2 = 20
$ws.Cells.Item($newRow, 18).Value  = "Hortaliza"
